# Add a "notes" column (R) to the student progress report.
# Each student's accumulated notes are rendered as a single cell whose
# individual entries are separated by carriage returns, most-recent note
# first (i.e. notes are shown in reverse chronological order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R1").Value = "notes"
$ws.Range("R2").Value = ""
$ws.Range("R3").Value = "Note #3: Note 3 (2021-05-04)`rthird note!`rNote #2: Note 2 (2021-04-18)`rsecond note!`rNote #1: Note 1 (2021-04-01)`rfirst note!`r"
$ws.Range("R4").Value = ""
$ws.Range("R5").Value = ""
$ws.Range("R6").Value = ""
$ws.Range("R7").Value = ""
